$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Avans (advance) values for July were moved to June, so clear column H
# for the rows that still had a carried-over advance amount.
$ws.Range("H3").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("H10").ClearContents()

# Update the active selection as recorded by the author when they saved.
$ws.Range("H18").Select()
